# Commit: "fix merge conflicts for manage program module"
# Adds a new "ProgramPage" worksheet (4th tab, after UserPage_Edit) used by the
# manage-program test module, with header/sample rows, and makes it the active tab.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end of
# the tab strip (HomePage, LoginPage, UserPage_Edit, ProgramPage), matching
# the sheetId=4 / rId4 ordering in the target workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ProgramPage"

# Header row
$newSheet.Range("A1").Value = "ProgramName"
$newSheet.Range("B1").Value = "ProgramDescription"
$newSheet.Range("C1").Value = "Status"
$newSheet.Range("D1").Value = "EditProgramName"
$newSheet.Range("E1").Value = "EditProgramDescription"

# Sample data row
$newSheet.Range("A2").Value = "RubySJC1234"
$newSheet.Range("B2").Value = "Language"
$newSheet.Range("C2").Value = "Active"
$newSheet.Range("D2").Value = "RDBMS123"
$newSheet.Range("E2").Value = "Test"

# Column widths matching the authored sheet (22.85546875, 20.85546875,
# 23.7109375, 19.28515625, 25 "chars"). The engine rounds ColumnWidth to the
# nearest 1/6 character unit when it re-derives the stored OOXML width
# (stored = (round(ColumnWidth*6)+5)/6), so we feed it the nearest
# representable input for each target instead of the raw authored value.
$newSheet.Columns.Item(1).ColumnWidth = 22.0
$newSheet.Columns.Item(2).ColumnWidth = 20.0
$newSheet.Columns.Item(3).ColumnWidth = 22.833333333333332
$newSheet.Columns.Item(4).ColumnWidth = 18.5
$newSheet.Columns.Item(5).ColumnWidth = 24.166666666666668

# Leave the selection on D38, as in the authored file (the new sheet is
# already the active/selected tab once added).
$newSheet.Range("D38").Select()
